{"js": "// Map of old division expressions -> new division expressions\nconst replacements = [\n  [\"18\u00f76=\", \"56\u00f73=\"],\n  [\"25\u00f77=\", \"41\u00f76=\"],\n  [\"17\u00f72=\", \"49\u00f75=\"],\n  [\"30\u00f73=\", \"85\u00f74=\"],\n  [\"37\u00f77=\", \"18\u00f77=\"],\n  [\"58\u00f72=\", \"25\u00f78=\"],\n  [\"17\u00f78=\", \"92\u00f78=\"],\n  [\"13\u00f74=\", \"21\u00f76=\"],\n  [\"42\u00f74=\", \"12\u00f73=\"],\n  [\"91\u00f72=\", \"61\u00f78=\"],\n  [\"72\u00f73=\", \"22\u00f74=\"],\n  [\"23\u00f79=\", \"93\u00f74=\"],\n  [\"90\u00f72=\", \"48\u00f77=\"],\n  [\"24\u00f76=\", \"80\u00f75=\"],\n  [\"98\u00f78=\", \"54\u00f78=\"],\n  [\"68\u00f74=\", \"55\u00f77=\"],\n  [\"47\u00f74=\", \"46\u00f73=\"],\n  [\"46\u00f72=\", \"54\u00f77=\"],\n  [\"35\u00f74=\", \"56\u00f79=\"],\n  [\"61\u00f75=\", \"94\u00f78=\"],\n  [\"71\u00f79=\", \"67\u00f77=\"],\n  [\"21\u00f74=\", \"34\u00f73=\"],\n  [\"59\u00f73=\", \"96\u00f77=\"],\n  [\"12\u00f74=\", \"16\u00f76=\"],\n  [\"95\u00f75=\", \"90\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"18\u00f76=\", \"56\u00f73=\"),\n    @(\"25\u00f77=\", \"41\u00f76=\"),\n    @(\"17\u00f72=\", \"49\u00f75=\"),\n    @(\"30\u00f73=\", \"85\u00f74=\"),\n    @(\"37\u00f77=\", \"18\u00f77=\"),\n    @(\"58\u00f72=\", \"25\u00f78=\"),\n    @(\"17\u00f78=\", \"92\u00f78=\"),\n    @(\"13\u00f74=\", \"21\u00f76=\"),\n    @(\"42\u00f74=\", \"12\u00f73=\"),\n    @(\"91\u00f72=\", \"61\u00f78=\"),\n    @(\"72\u00f73=\", \"22\u00f74=\"),\n    @(\"23\u00f79=\", \"93\u00f74=\"),\n    @(\"90\u00f72=\", \"48\u00f77=\"),\n    @(\"24\u00f76=\", \"80\u00f75=\"),\n    @(\"98\u00f78=\", \"54\u00f78=\"),\n    @(\"68\u00f74=\", \"55\u00f77=\"),\n    @(\"47\u00f74=\", \"46\u00f73=\"),\n    @(\"46\u00f72=\", \"54\u00f77=\"),\n    @(\"35\u00f74=\", \"56\u00f79=\"),\n    @(\"61\u00f75=\", \"94\u00f78=\"),\n    @(\"71\u00f79=\", \"67\u00f77=\"),\n    @(\"21\u00f74=\", \"34\u00f73=\"),\n    @(\"59\u00f73=\", \"96\u00f77=\"),\n    @(\"12\u00f74=\", \"16\u00f76=\"),\n    @(\"95\u00f75=\", \"90\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
